$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-04-05 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-04-06 Sunday", 2) | Out-Null
$d.Content.Find.Execute("79×73=5767", $true, $false, $false, $false, $false, $true, 1, $false, "89×58=5162", 2) | Out-Null
$d.Content.Find.Execute("17×70=1190", $true, $false, $false, $false, $false, $true, 1, $false, "84×41=3444", 2) | Out-Null
$d.Content.Find.Execute("92×17=1564", $true, $false, $false, $false, $false, $true, 1, $false, "54×55=2970", 2) | Out-Null
$d.Content.Find.Execute("32×86=2752", $true, $false, $false, $false, $false, $true, 1, $false, "18×77=1386", 2) | Out-Null
$d.Content.Find.Execute("73×65=4745", $true, $false, $false, $false, $false, $true, 1, $false, "31×21=651", 2) | Out-Null
$d.Content.Find.Execute("16×89=1424", $true, $false, $false, $false, $false, $true, 1, $false, "90×42=3780", 2) | Out-Null
$d.Content.Find.Execute("52×71=3692", $true, $false, $false, $false, $false, $true, 1, $false, "83×59=4897", 2) | Out-Null
$d.Content.Find.Execute("31×92=2852", $true, $false, $false, $false, $false, $true, 1, $false, "15×49=735", 2) | Out-Null
$d.Content.Find.Execute("70×38=2660", $true, $false, $false, $false, $false, $true, 1, $false, "25×90=2250", 2) | Out-Null
$d.Content.Find.Execute("97×92=8924", $true, $false, $false, $false, $false, $true, 1, $false, "24×46=1104", 2) | Out-Null
$d.Content.Find.Execute("82×48=3936", $true, $false, $false, $false, $false, $true, 1, $false, "83×25=2075", 2) | Out-Null
$d.Content.Find.Execute("42×50=2100", $true, $false, $false, $false, $false, $true, 1, $false, "87×31=2697", 2) | Out-Null
$d.Content.Find.Execute("11×62=682", $true, $false, $false, $false, $false, $true, 1, $false, "11×54=594", 2) | Out-Null
$d.Content.Find.Execute("13×32=416", $true, $false, $false, $false, $false, $true, 1, $false, "47×95=4465", 2) | Out-Null
$d.Content.Find.Execute("78×76=5928", $true, $false, $false, $false, $false, $true, 1, $false, "16×99=1584", 2) | Out-Null
$d.Content.Find.Execute("84×54=4536", $true, $false, $false, $false, $false, $true, 1, $false, "59×35=2065", 2) | Out-Null
$d.Content.Find.Execute("25×42=1050", $true, $false, $false, $false, $false, $true, 1, $false, "24×93=2232", 2) | Out-Null
$d.Content.Find.Execute("79×49=3871", $true, $false, $false, $false, $false, $true, 1, $false, "17×35=595", 2) | Out-Null
$d.Content.Find.Execute("21×94=1974", $true, $false, $false, $false, $false, $true, 1, $false, "43×67=2881", 2) | Out-Null
$d.Content.Find.Execute("74×20=1480", $true, $false, $false, $false, $false, $true, 1, $false, "33×58=1914", 2) | Out-Null
$d.Content.Find.Execute("14×98=1372", $true, $false, $false, $false, $false, $true, 1, $false, "40×22=880", 2) | Out-Null
$d.Content.Find.Execute("65×49=3185", $true, $false, $false, $false, $false, $true, 1, $false, "46×20=920", 2) | Out-Null
$d.Content.Find.Execute("17×27=459", $true, $false, $false, $false, $false, $true, 1, $false, "58×98=5684", 2) | Out-Null
$d.Content.Find.Execute("92×43=3956", $true, $false, $false, $false, $false, $true, 1, $false, "51×80=4080", 2) | Out-Null
$d.Content.Find.Execute("94×51=4794", $true, $false, $false, $false, $false, $true, 1, $false, "22×51=1122", 2) | Out-Null
